$wb = $excel.ActiveWorkbook

# Sheet ALC row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1213.9286
$ws.Range("J112").Value = 1219.6
$ws.Range("L112").Value = 3658.8
$ws.Range("N112").Value = -5874.799999999999

# Sheet ALC row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1238.1724
$ws.Range("J129").Value = 1238.1724
$ws.Range("L129").Value = 3714.5172
$ws.Range("N129").Value = -13714.5172

# Sheet ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1098.2
$ws.Range("I137").Value = 995.15625
$ws.Range("J137").Value = 1510.375
$ws.Range("K137").Value = 2985.46875
$ws.Range("L137").Value = 4531.125
$ws.Range("M137").Value = -435.46875
$ws.Range("N137").Value = -9631.125

# Sheet ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2366.13
$ws.Range("I138").Value = 1174.0896
$ws.Range("J138").Value = 4786.3335
$ws.Range("K138").Value = 3522.2688
$ws.Range("L138").Value = 14359.0005
$ws.Range("M138").Value = 1617.7312
$ws.Range("N138").Value = -24639.0005

# Sheet ARM row 23
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 53336.89
$ws.Range("J23").Value = 45002.332
$ws.Range("L23").Value = 45002.332
$ws.Range("N23").Value = -45520.332

# Sheet ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4509.737
$ws.Range("I32").Value = 3223.5059
$ws.Range("J32").Value = 18497.5
$ws.Range("K32").Value = 3223.5059
$ws.Range("L32").Value = 18497.5
$ws.Range("M32").Value = -2936.5059
$ws.Range("N32").Value = -19071.5

# Sheet ARM row 37
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 12111.111
$ws.Range("I37").Value = 12111.111
$ws.Range("K37").Value = 12111.111
$ws.Range("M37").Value = -11838.111

# Sheet ARM row 44
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 30841.666
$ws.Range("J44").Value = 30841.666
$ws.Range("L44").Value = 30841.666
$ws.Range("N44").Value = -31817.666

# Sheet ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1217.9286
$ws.Range("I45").Value = 1088.5
$ws.Range("K45").Value = 1088.5
$ws.Range("M45").Value = -711.5

# Sheet ARM row 55
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 9948
$ws.Range("I55").Value = 9948
$ws.Range("K55").Value = 9948
$ws.Range("M55").Value = -9633

# Sheet ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 777.5599999999999
$ws.Range("I61").Value = 732.65216
$ws.Range("K61").Value = 732.65216
$ws.Range("M61").Value = -520.65216

# Sheet ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3166.8235
$ws.Range("J74").Value = 1280.5714
$ws.Range("L74").Value = 1280.5714
$ws.Range("N74").Value = -3028.5714

# Sheet ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 3166.8235
$ws.Range("J77").Value = 1280.5714
$ws.Range("L77").Value = 6402.857
$ws.Range("N77").Value = -15138.857

# Sheet ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 40354
$ws.Range("I110").Value = 56735.777
$ws.Range("J110").Value = 3495
$ws.Range("K110").Value = 56735.777
$ws.Range("L110").Value = 3495
$ws.Range("M110").Value = -54690.777
$ws.Range("N110").Value = -7585

# Sheet ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2402.4546
$ws.Range("I122").Value = 1289.5714
$ws.Range("J122").Value = 4350
$ws.Range("K122").Value = 3868.7142
$ws.Range("L122").Value = 13050
$ws.Range("M122").Value = -1418.7142
$ws.Range("N122").Value = -17950

# Sheet ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 777.5599999999999
$ws.Range("I136").Value = 732.65216
$ws.Range("K136").Value = 2197.95648
$ws.Range("M136").Value = 352.0435200000002

# Sheet CRP row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3173
$ws.Range("I16").Value = 3319.625
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 3319.625
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -3032.625
$ws.Range("N16").Value = -2574

# Sheet CRP row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 275
$ws.Range("I22").Value = 175
$ws.Range("J22").Value = 375
$ws.Range("K22").Value = 175
$ws.Range("L22").Value = 375
$ws.Range("M22").Value = 175
$ws.Range("N22").Value = -1075

# Sheet CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3030.7917
$ws.Range("I31").Value = 2203.9443
$ws.Range("J31").Value = 3526.9
$ws.Range("K31").Value = 2203.9443
$ws.Range("L31").Value = 3526.9
$ws.Range("M31").Value = -1908.9443
$ws.Range("N31").Value = -4116.9

# Sheet CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3030.7917
$ws.Range("I34").Value = 2203.9443
$ws.Range("J34").Value = 3526.9
$ws.Range("K34").Value = 2203.9443
$ws.Range("L34").Value = 3526.9
$ws.Range("M34").Value = -2001.9443
$ws.Range("N34").Value = -3930.9

# Sheet CRP row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 3173
$ws.Range("I113").Value = 3319.625
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 3319.625
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = -1149.625
$ws.Range("N113").Value = -6340

# Sheet CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 3735
$ws.Range("I131").Value = 375
$ws.Range("J131").Value = 6741.316
$ws.Range("K131").Value = 1125
$ws.Range("L131").Value = 20223.948
$ws.Range("M131").Value = 3915
$ws.Range("N131").Value = -30303.948

# Sheet CUL row 140
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 2538.6667
$ws.Range("I140").Value = 2531.6667
$ws.Range("J140").Value = 2543.3333
$ws.Range("K140").Value = 7595.000100000001
$ws.Range("L140").Value = 7629.999899999999
$ws.Range("M140").Value = -2415.000100000001
$ws.Range("N140").Value = -17989.9999

# Sheet GSM row 2
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 56.75
$ws.Range("I2").Value = 63.5
$ws.Range("J2").Value = 50
$ws.Range("K2").Value = 63.5
$ws.Range("L2").Value = 50
$ws.Range("M2").Value = 49.5
$ws.Range("N2").Value = -276

# Sheet GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3800
$ws.Range("I102").Value = 800
$ws.Range("J102").Value = 4800
$ws.Range("K102").Value = 800
$ws.Range("L102").Value = 4800
$ws.Range("M102").Value = 822
$ws.Range("N102").Value = -8044

# Sheet GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2193.3333
$ws.Range("I122").Value = 1746
$ws.Range("K122").Value = 5238
$ws.Range("M122").Value = -2788

# Sheet LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 586.8461
$ws.Range("I22").Value = 628.1667
$ws.Range("J22").Value = 551.4286
$ws.Range("K22").Value = 628.1667
$ws.Range("L22").Value = 551.4286
$ws.Range("M22").Value = -333.1667
$ws.Range("N22").Value = -1141.4286

# Sheet LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 586.8461
$ws.Range("I27").Value = 628.1667
$ws.Range("J27").Value = 551.4286
$ws.Range("K27").Value = 628.1667
$ws.Range("L27").Value = 551.4286
$ws.Range("M27").Value = -521.1667
$ws.Range("N27").Value = -765.4286

# Sheet LTW row 34
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value = 50000
$ws.Range("J34").Value = 50000
$ws.Range("L34").Value = 50000
$ws.Range("N34").Value = -50344

# Sheet LTW row 55
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 682.75
$ws.Range("I55").Value = 400.2
$ws.Range("J55").Value = 884.5714
$ws.Range("K55").Value = 400.2
$ws.Range("L55").Value = 884.5714
$ws.Range("M55").Value = -227.2
$ws.Range("N55").Value = -1230.5714

# Sheet LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4948.0244
$ws.Range("I132").Value = 4733.3276
$ws.Range("J132").Value = 5466.875
$ws.Range("K132").Value = 14199.9828
$ws.Range("L132").Value = 16400.625
$ws.Range("M132").Value = -11669.9828
$ws.Range("N132").Value = -21460.625

# Sheet WVR row 14
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 3500
$ws.Range("J14").Value = 3500
$ws.Range("L14").Value = 3500
$ws.Range("N14").Value = -3836

# Sheet WVR row 45
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 20626
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 20626
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 20626
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -21608

# Sheet WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1442
$ws.Range("I132").Value = 1442
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4326
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1796
$ws.Range("N132").ClearContents()

# Sheet WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1298.0546
$ws.Range("I136").Value = 470.35294
$ws.Range("J136").Value = 11851.25
$ws.Range("K136").Value = 1411.05882
$ws.Range("L136").Value = 35553.75
$ws.Range("M136").Value = 1138.94118
$ws.Range("N136").Value = -40653.75
